# Edit script for assets/disciplinas/8800012.xlsx
#
# The commit reshuffles the content of rows 10 and 13-23 of the single
# worksheet: most "label" cells in column A shift up by one row while the
# associated value cells in columns B/C are resupplied from a different
# (nearby) row, one brand new value ("Semestral") is introduced, and the
# final row (24) -- which only held a duplicate "8800011..." requirement
# note -- is removed so the sheet shrinks from A1:C24 to A1:C23.
#
# To apply this safely we first snapshot every value we still need (from
# column A and B) before overwriting anything, then write the new layout
# back using those snapshots. The "01/01/2017" value is moved with
# .Copy() instead of being retyped, so Excel keeps treating it as text
# instead of silently re-parsing it as a date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- snapshot values that are needed later, before they get overwritten ----
$A10 = $ws.Cells.Item(10, 1).Value()
$A14 = $ws.Cells.Item(14, 1).Value()
$A15 = $ws.Cells.Item(15, 1).Value()
$A16 = $ws.Cells.Item(16, 1).Value()
$A17 = $ws.Cells.Item(17, 1).Value()
$A18 = $ws.Cells.Item(18, 1).Value()
$A19 = $ws.Cells.Item(19, 1).Value()
$A20 = $ws.Cells.Item(20, 1).Value()
$A21 = $ws.Cells.Item(21, 1).Value()
$A22 = $ws.Cells.Item(22, 1).Value()
$A23 = $ws.Cells.Item(23, 1).Value()

$B13 = $ws.Cells.Item(13, 2).Value()
$B15 = $ws.Cells.Item(15, 2).Value()
$B17 = $ws.Cells.Item(17, 2).Value()
$B19 = $ws.Cells.Item(19, 2).Value()
$B20 = $ws.Cells.Item(20, 2).Value()
$B21 = $ws.Cells.Item(21, 2).Value()
$B24 = $ws.Cells.Item(24, 2).Value()

# ---- row 10: Objetivos: now shows the docent name instead of the goals ----
$ws.Cells.Item(10, 2).Value() = $B13
$ws.Cells.Item(10, 3).Value() = $B13

# ---- row 13: gains the "Programa resumido:" label and a new "Semestral" value ----
$ws.Cells.Item(13, 1).Value() = $A14
$ws.Cells.Item(13, 2).Value() = "Semestral"
$ws.Cells.Item(13, 3).Value() = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---- row 14: becomes "Short syllabus:" with its (already adjacent) English text ----
$ws.Cells.Item(14, 1).Value() = $A15
$ws.Cells.Item(14, 2).Value() = $B15
$ws.Cells.Item(14, 3).Value() = $B15

# ---- row 15: becomes "Programa:" paired with the "01/01/2017" text value ----
$ws.Cells.Item(15, 1).Value() = $A16
$ws.Cells.Item(8, 2).Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(8, 3).Copy($ws.Cells.Item(15, 3))
$ws.Rows.Item(15).RowHeight = 120

# ---- row 16: becomes "Syllabus:" with its English text ----
$ws.Cells.Item(16, 1).Value() = $A17
$ws.Cells.Item(16, 2).Value() = $B17
$ws.Cells.Item(16, 3).Value() = $B17
$ws.Rows.Item(16).RowHeight = 120

# ---- row 17: becomes "Avaliação:" (label only, default row height) ----
$ws.Cells.Item(17, 1).Value() = $A18
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).ClearContents()
$ws.Rows.Item(17).EntireRow.AutoFit()

# ---- row 18: becomes "Método:" paired with the docent name ----
$ws.Cells.Item(18, 1).Value() = $A19
$ws.Cells.Item(18, 2).Value() = $B13
$ws.Cells.Item(18, 3).Value() = $B13
$ws.Rows.Item(18).RowHeight = 60

# ---- row 19: becomes "Critério:" with the semester-program text ----
$ws.Cells.Item(19, 1).Value() = $A20
$ws.Cells.Item(19, 2).Value() = $B19
$ws.Cells.Item(19, 3).Value() = $B19

# ---- row 20: becomes "Norma de recuperação:" with the evaluation-criteria text ----
$ws.Cells.Item(20, 1).Value() = $A21
$ws.Cells.Item(20, 2).Value() = $B20
$ws.Cells.Item(20, 3).Value() = $B20

# ---- row 21: becomes "Bibliografia:" with "Não tem" ----
$ws.Cells.Item(21, 1).Value() = $A22
$ws.Cells.Item(21, 2).Value() = $B21
$ws.Cells.Item(21, 3).Value() = $B21
$ws.Rows.Item(21).RowHeight = 120

# ---- row 22: becomes "Requisitos:" (label only, default row height) ----
$ws.Cells.Item(22, 1).Value() = $A23
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).ClearContents()
$ws.Rows.Item(22).EntireRow.AutoFit()

# ---- row 23: loses its label, keeps the "8800011..." requirement note ----
$ws.Cells.Item(23, 1).ClearContents()
$ws.Cells.Item(23, 2).Value() = $B24
$ws.Cells.Item(23, 3).Value() = $B24
$ws.Rows.Item(23).RowHeight = 30

# ---- the old trailing row 24 is now redundant; delete it ----
$ws.Rows.Item(24).Delete()
